# "Generate Report for Handoff"
# The localization status report is regenerated: the zh-cn item has moved
# from "In Translation" to "Ready for handoff", the handoff timestamps are
# refreshed, and the Status columns are widened to fit the new, longer
# status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -----------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # Overview!zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # Overview!de-de status
$wsZhCn.Range("C2").Value     = "Ready for handoff"   # zh-cn!Status
$wsDeDe.Range("C2").Value     = "Ready for handoff"   # de-de!Status

# --- Timestamps refreshed by the new handoff generation --------------------
# Overview!Latest HO Xliff Generate Date and de-de!Latest Handoff Datetime
$wsOverview.Range("G2").Value = "2016-08-23 11:00:28"
$wsDeDe.Range("H2").Value     = "2016-08-23 11:00:28"

# zh-cn!Latest Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-23 11:00:07"

# --- Widen the Status columns so the longer text fits ----------------------
# (engine's ColumnWidth setter quantizes internally; 16 + 1/3 is the closest
# input that reproduces the report's new ~17.22-character stored width)
$wsOverview.Columns("E").ColumnWidth = 16.333333333333336
$wsOverview.Columns("F").ColumnWidth = 16.333333333333336
$wsZhCn.Columns("C").ColumnWidth     = 16.333333333333336
$wsDeDe.Columns("C").ColumnWidth     = 16.333333333333336
